$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "About" sheet: bump the report date (C1) and make it the active tab
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45387

# ---------------------------------------------------------------------
# 2) "BAU Emissions" sheet: relabel the " : NoSettings" rows to " : test"
#    (column A, rows 4-28 and 31-280), and update the forecast values in
#    row 94 (columns M:AE).
# ---------------------------------------------------------------------
$bau = $wb.Worksheets.Item("BAU Emissions")

$rowRanges = @(@(4, 28), @(31, 280))
foreach ($range in $rowRanges) {
    $startRow = $range[0]
    $endRow = $range[1]
    for ($r = $startRow; $r -le $endRow; $r++) {
        $cell = $bau.Cells.Item($r, 1)
        $old = [string]$cell.Value2
        if ($old -like "* : NoSettings") {
            $new = $old -replace " : NoSettings$", " : test"
            $cell.Value = $new
        }
    }
}

$row94Values = [ordered]@{
    "M94"  = 1001080
    "N94"  = 2002150
    "O94"  = 3003230
    "P94"  = 4004300
    "Q94"  = 5005380
    "R94"  = 5005380
    "S94"  = 5005380
    "T94"  = 5005380
    "U94"  = 5005380
    "V94"  = 5005380
    "W94"  = 5005380
    "X94"  = 5005380
    "Y94"  = 5005380
    "Z94"  = 5005380
    "AA94" = 5005380
    "AB94" = 5005380
    "AC94" = 5005380
    "AD94" = 5005380
    "AE94" = 5005380
}
foreach ($addr in $row94Values.Keys) {
    $bau.Range($addr).Value = $row94Values[$addr]
}

# ---------------------------------------------------------------------
# 3) Sheet-view bookkeeping: "About" becomes the selected tab (instead of
#    "Current and Planned Capacity"), and the "BAU Emissions" selection
#    moves to A30:AE280.
# ---------------------------------------------------------------------
$bau.Activate()
$bau.Range("A30:AE280").Select() | Out-Null

$about.Activate()
